$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 403.83334
$ws.Range("J2").Value = 603.8
$ws.Range("L2").Value = 603.8
$ws.Range("N2").Value = -829.8
$ws.Range("H12").Value = 7398.2856
$ws.Range("I12").Value = 8599.25
$ws.Range("J12").Value = 192.5
$ws.Range("K12").Value = 8599.25
$ws.Range("L12").Value = 192.5
$ws.Range("M12").Value = -8429.25
$ws.Range("N12").Value = -532.5
$ws.Range("H13").Value = 29999
$ws.Range("J13").Value = 29999
$ws.Range("L13").Value = 29999
$ws.Range("N13").Value = -30337
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("H32").Value = 564.5
$ws.Range("I32").Value = 556
$ws.Range("K32").Value = 556
$ws.Range("M32").Value = -230
$ws.Range("H53").Value = 374.8
$ws.Range("I53").Value = 57.833332
$ws.Range("K53").Value = 57.833332
$ws.Range("M53").Value = 579.166668
$ws.Range("H55").Value = 109
$ws.Range("J55").Value = 88
$ws.Range("L55").Value = 88
$ws.Range("N55").Value = -516
$ws.Range("H64").Value = 16670932
$ws.Range("I64").Value = 20837456
$ws.Range("J64").Value = 4833.3335
$ws.Range("K64").Value = 20837456
$ws.Range("L64").Value = 4833.3335
$ws.Range("M64").Value = -20837208
$ws.Range("N64").Value = -5329.3335
$ws.Range("H67").Value = 16670932
$ws.Range("I67").Value = 20837456
$ws.Range("J67").Value = 4833.3335
$ws.Range("K67").Value = 20837456
$ws.Range("L67").Value = 4833.3335
$ws.Range("M67").Value = -20836598
$ws.Range("N67").Value = -6549.3335
$ws.Range("H74").Value = 62508000
$ws.Range("I74").Value = 125012000
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 125012000
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -125011064
$ws.Range("N74").Value = -5872
$ws.Range("H77").Value = 62508000
$ws.Range("I77").Value = 125012000
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 625060000
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -625055320
$ws.Range("N77").Value = -29360
$ws.Range("H80").Value = 2211
$ws.Range("I80").Value = 3064.3333
$ws.Range("J80").Value = 504.33334
$ws.Range("K80").Value = 9192.999899999999
$ws.Range("L80").Value = 1513.00002
$ws.Range("M80").Value = -8194.999899999999
$ws.Range("N80").Value = -3509.00002
$ws.Range("H83").Value = 2211
$ws.Range("I83").Value = 3064.3333
$ws.Range("J83").Value = 504.33334
$ws.Range("K83").Value = 27578.9997
$ws.Range("L83").Value = 4539.00006
$ws.Range("M83").Value = -22586.9997
$ws.Range("N83").Value = -14523.00006
$ws.Range("H88").Value = 52734090
$ws.Range("I88").Value = 333333340
$ws.Range("J88").Value = 5967546.5
$ws.Range("K88").Value = 333333340
$ws.Range("L88").Value = 5967546.5
$ws.Range("M88").Value = -333332934
$ws.Range("N88").Value = -5968358.5
$ws.Range("H91").Value = 52734090
$ws.Range("I91").Value = 333333340
$ws.Range("J91").Value = 5967546.5
$ws.Range("K91").Value = 333333340
$ws.Range("L91").Value = 5967546.5
$ws.Range("M91").Value = -333331936
$ws.Range("N91").Value = -5970354.5
$ws.Range("H94").Value = 1784.2858
$ws.Range("I94").Value = 1784.2858
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1784.2858
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1333.2858
$ws.Range("H111").Value = 201539.8
$ws.Range("I111").Value = 1100
$ws.Range("K111").Value = 3300
$ws.Range("M111").Value = -233
$ws.Range("H113").Value = 3128.2
$ws.Range("I113").Value = 3128.2
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3128.2
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 125.8000000000002
$ws.Range("H116").Value = 7555.125
$ws.Range("I116").Value = 7648.5654
$ws.Range("J116").Value = 5406
$ws.Range("K116").Value = 7648.5654
$ws.Range("L116").Value = 5406
$ws.Range("M116").Value = -4206.5654
$ws.Range("N116").Value = -12290
$ws.Range("H129").Value = 9398.546
$ws.Range("I129").Value = 1738.7
$ws.Range("K129").Value = 5216.1
$ws.Range("M129").Value = -216.1000000000004
$ws.Range("H132").Value = 3529.7273
$ws.Range("I132").Value = 3682.7
$ws.Range("K132").Value = 11048.1
$ws.Range("M132").Value = -8518.099999999999
$ws.Range("H135").Value = 35714630
$ws.Range("I135").Value = 35714630
$ws.Range("K135").Value = 321431670
$ws.Range("M135").Value = -321429135
$ws.Range("H137").Value = 2315.9
$ws.Range("I137").Value = 1456.4615
$ws.Range("K137").Value = 4369.3845
$ws.Range("M137").Value = -1819.3845
$ws.Range("H138").Value = 3019.2903
$ws.Range("I138").Value = 2164.6924
$ws.Range("J138").Value = 3636.5
$ws.Range("K138").Value = 6494.0772
$ws.Range("L138").Value = 10909.5
$ws.Range("M138").Value = -1354.0772
$ws.Range("N138").Value = -21189.5
$ws.Range("H141").Value = 2926.182
$ws.Range("I141").Value = 2926.182
$ws.Range("K141").Value = 8778.545999999998
$ws.Range("M141").Value = -3598.545999999998
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("N113").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1228772
$ws.Range("J2").Value = 2473.5
$ws.Range("L2").Value = 2473.5
$ws.Range("N2").Value = -2699.5
$ws.Range("H5").Value = 2010
$ws.Range("I5").Value = 1101
$ws.Range("J5").Value = 3525
$ws.Range("K5").Value = 1101
$ws.Range("L5").Value = 3525
$ws.Range("M5").Value = -989
$ws.Range("N5").Value = -3749
$ws.Range("H32").Value = 7044.96
$ws.Range("I32").Value = 5509.7393
$ws.Range("J32").Value = 24700
$ws.Range("K32").Value = 5509.7393
$ws.Range("L32").Value = 24700
$ws.Range("M32").Value = -5222.7393
$ws.Range("N32").Value = -25274
$ws.Range("H61").Value = 47621310
$ws.Range("I61").Value = 71429896
$ws.Range("K61").Value = 71429896
$ws.Range("M61").Value = -71429684
$ws.Range("H74").Value = 43480764
$ws.Range("I74").Value = 47620836
$ws.Range("K74").Value = 47620836
$ws.Range("M74").Value = -47619962
$ws.Range("H77").Value = 43480764
$ws.Range("I77").Value = 47620836
$ws.Range("K77").Value = 238104180
$ws.Range("M77").Value = -238099812
$ws.Range("H88").Value = 79448.53999999999
$ws.Range("I88").Value = 127042.875
$ws.Range("J88").Value = 3297.6
$ws.Range("K88").Value = 127042.875
$ws.Range("L88").Value = 3297.6
$ws.Range("M88").Value = -126636.875
$ws.Range("N88").Value = -4109.6
$ws.Range("H91").Value = 79448.53999999999
$ws.Range("I91").Value = 127042.875
$ws.Range("J91").Value = 3297.6
$ws.Range("K91").Value = 127042.875
$ws.Range("L91").Value = 3297.6
$ws.Range("M91").Value = -125638.875
$ws.Range("N91").Value = -6105.6
$ws.Range("H97").Value = 611.35297
$ws.Range("I97").Value = 769.4783
$ws.Range("K97").Value = 769.4783
$ws.Range("M97").Value = -273.4783
$ws.Range("H105").Value = 44437.5
$ws.Range("J105").Value = 44437.5
$ws.Range("L105").Value = 44437.5
$ws.Range("N105").Value = -51425.5
$ws.Range("H116").Value = 1228772
$ws.Range("J116").Value = 2473.5
$ws.Range("L116").Value = 2473.5
$ws.Range("N116").Value = -7061.5
$ws.Range("H132").Value = 5266140.5
$ws.Range("I132").Value = 7145033.5
$ws.Range("J132").Value = 5239.2
$ws.Range("K132").Value = 21435100.5
$ws.Range("L132").Value = 15717.6
$ws.Range("M132").Value = -21432570.5
$ws.Range("N132").Value = -20777.6
$ws.Range("H136").Value = 47621310
$ws.Range("I136").Value = 71429896
$ws.Range("K136").Value = 214289688
$ws.Range("M136").Value = -214287138

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1228772
$ws.Range("J3").Value = 2473.5
$ws.Range("L3").Value = 2473.5
$ws.Range("N3").Value = -2701.5
$ws.Range("H4").Value = 2010
$ws.Range("I4").Value = 1101
$ws.Range("J4").Value = 3525
$ws.Range("K4").Value = 1101
$ws.Range("L4").Value = 3525
$ws.Range("M4").Value = -986
$ws.Range("N4").Value = -3755
$ws.Range("H24").Value = 11253.75
$ws.Range("I24").Value = 8338.666999999999
$ws.Range("K24").Value = 8338.666999999999
$ws.Range("M24").Value = -8103.666999999999
$ws.Range("H75").Value = 11523.857
$ws.Range("J75").Value = 17736
$ws.Range("L75").Value = 17736
$ws.Range("N75").Value = -19608
$ws.Range("H78").Value = 11523.857
$ws.Range("J78").Value = 17736
$ws.Range("L78").Value = 53208
$ws.Range("N78").Value = -62568
$ws.Range("H86").Value = 2201.8572
$ws.Range("I86").Value = 2294.3076
$ws.Range("K86").Value = 2294.3076
$ws.Range("M86").Value = -1171.3076
$ws.Range("H89").Value = 2201.8572
$ws.Range("I89").Value = 2294.3076
$ws.Range("K89").Value = 11471.538
$ws.Range("M89").Value = -5855.538
$ws.Range("H99").Value = 2140.7856
$ws.Range("I99").Value = 1891
$ws.Range("J99").Value = 2240.7
$ws.Range("K99").Value = 1891
$ws.Range("L99").Value = 2240.7
$ws.Range("M99").Value = -393
$ws.Range("N99").Value = -5236.7
$ws.Range("H105").Value = 2851.818
$ws.Range("I105").Value = 1895
$ws.Range("K105").Value = 1895
$ws.Range("M105").Value = -148
$ws.Range("H107").Value = 155396.31
$ws.Range("I107").Value = 1832.2727
$ws.Range("J107").Value = 999998.5
$ws.Range("K107").Value = 1832.2727
$ws.Range("L107").Value = 999998.5
$ws.Range("M107").Value = 87.72730000000001
$ws.Range("N107").Value = -1003838.5
$ws.Range("H108").Value = 69998.664
$ws.Range("J108").Value = 69998.664
$ws.Range("L108").Value = 69998.664
$ws.Range("N108").Value = -77678.664
$ws.Range("H134").Value = 20839238
$ws.Range("I134").Value = 21745074
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 65235222
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -65232687
$ws.Range("N134").Value = -20070

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1036.8572
$ws.Range("I7").Value = 1339.1333
$ws.Range("J7").Value = 281.16666
$ws.Range("K7").Value = 1339.1333
$ws.Range("L7").Value = 281.16666
$ws.Range("M7").Value = -1226.1333
$ws.Range("N7").Value = -507.16666
$ws.Range("H31").Value = 3715.4546
$ws.Range("I31").Value = 3661.111
$ws.Range("K31").Value = 3661.111
$ws.Range("M31").Value = -3366.111
$ws.Range("H34").Value = 3715.4546
$ws.Range("I34").Value = 3661.111
$ws.Range("K34").Value = 3661.111
$ws.Range("M34").Value = -3459.111
$ws.Range("H51").Value = 44998.555
$ws.Range("J51").Value = 44998.555
$ws.Range("L51").Value = 44998.555
$ws.Range("N51").Value = -46470.555
$ws.Range("H58").Value = 13518035
$ws.Range("J58").Value = 2792.6667
$ws.Range("L58").Value = 2792.6667
$ws.Range("N58").Value = -3198.6667
$ws.Range("H60").Value = 35030.69
$ws.Range("J60").Value = 39999.91
$ws.Range("L60").Value = 39999.91
$ws.Range("N60").Value = -41021.91
$ws.Range("H61").Value = 44998.555
$ws.Range("J61").Value = 44998.555
$ws.Range("L61").Value = 44998.555
$ws.Range("N61").Value = -45694.555
$ws.Range("H74").Value = 48519
$ws.Range("J74").Value = 49938.832
$ws.Range("L74").Value = 49938.832
$ws.Range("N74").Value = -51686.832
$ws.Range("H77").Value = 48519
$ws.Range("J77").Value = 49938.832
$ws.Range("L77").Value = 149816.496
$ws.Range("N77").Value = -158552.496
$ws.Range("H99").Value = 4497.8335
$ws.Range("I99").Value = 4596.75
$ws.Range("J99").Value = 4300
$ws.Range("K99").Value = 4596.75
$ws.Range("L99").Value = 4300
$ws.Range("M99").Value = -3098.75
$ws.Range("N99").Value = -7296
$ws.Range("H126").Value = 4497.8335
$ws.Range("I126").Value = 4596.75
$ws.Range("J126").Value = 4300
$ws.Range("K126").Value = 13790.25
$ws.Range("L126").Value = 12900
$ws.Range("M126").Value = -11320.25
$ws.Range("N126").Value = -17840
$ws.Range("H132").Value = 125002470
$ws.Range("I132").Value = 250001460
$ws.Range("K132").Value = 750004380
$ws.Range("M132").Value = -750001850
$ws.Range("H134").Value = 14767601
$ws.Range("I134").Value = 15690389
$ws.Range("K134").Value = 47071167
$ws.Range("M134").Value = -47068632
$ws.Range("H136").Value = 13518035
$ws.Range("J136").Value = 2792.6667
$ws.Range("L136").Value = 8378.000100000001
$ws.Range("N136").Value = -13478.0001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 297.75
$ws.Range("I12").Value = 39.375
$ws.Range("J12").Value = 470
$ws.Range("K12").Value = 118.125
$ws.Range("L12").Value = 1410
$ws.Range("M12").Value = 54.875
$ws.Range("N12").Value = -1756
$ws.Range("H14").Value = 224
$ws.Range("I14").Value = 224
$ws.Range("K14").Value = 672
$ws.Range("M14").Value = -499
$ws.Range("H33").Value = 899.93335
$ws.Range("I33").Value = 474.5
$ws.Range("K33").Value = 2847
$ws.Range("M33").Value = -2564
$ws.Range("H107").Value = 1102.875
$ws.Range("J107").Value = 1376.1666
$ws.Range("L107").Value = 4128.4998
$ws.Range("N107").Value = -7968.4998
$ws.Range("H113").Value = 63134.312
$ws.Range("I113").Value = 143408.72
$ws.Range("K113").Value = 430226.16
$ws.Range("M113").Value = -428056.16
$ws.Range("H121").Value = 103902.5
$ws.Range("J121").Value = 6763.4
$ws.Range("L121").Value = 20290.2
$ws.Range("N121").Value = -22910.2
$ws.Range("H125").Value = 11000
$ws.Range("J125").Value = 10000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840
$ws.Range("H131").Value = 1742.5714
$ws.Range("I131").Value = 1260.3
$ws.Range("K131").Value = 3780.9
$ws.Range("M131").Value = 1259.1
$ws.Range("H132").Value = 2387.7896
$ws.Range("I132").Value = 2149.5
$ws.Range("J132").Value = 2415.8235
$ws.Range("K132").Value = 19345.5
$ws.Range("L132").Value = 21742.4115
$ws.Range("M132").Value = -16815.5
$ws.Range("N132").Value = -26802.4115
$ws.Range("H140").Value = 1585.75
$ws.Range("I140").Value = 1585.75
$ws.Range("K140").Value = 4757.25
$ws.Range("M140").Value = 422.75
$ws.Range("H141").Value = 3998.5
$ws.Range("I141").Value = 3998.5
$ws.Range("K141").Value = 11995.5
$ws.Range("M141").Value = -6815.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 176.72728
$ws.Range("I2").Value = 144.33333
$ws.Range("K2").Value = 144.33333
$ws.Range("M2").Value = -31.33332999999999
$ws.Range("H43").Value = 1665
$ws.Range("I43").Value = 1665
$ws.Range("K43").Value = 1665
$ws.Range("M43").Value = -1514
$ws.Range("H97").Value = 3282.923
$ws.Range("I97").Value = 3368.3
$ws.Range("J97").Value = 2998.3333
$ws.Range("K97").Value = 3368.3
$ws.Range("L97").Value = 2998.3333
$ws.Range("M97").Value = -2872.3
$ws.Range("N97").Value = -3990.3333
$ws.Range("H102").Value = 5914
$ws.Range("I102").Value = 4771.143
$ws.Range("K102").Value = 4771.143
$ws.Range("M102").Value = -3149.143
$ws.Range("H122").Value = 71052.94500000001
$ws.Range("I122").Value = 97693.46000000001
$ws.Range("J122").Value = 13331.833
$ws.Range("K122").Value = 293080.38
$ws.Range("L122").Value = 39995.499
$ws.Range("M122").Value = -290630.38
$ws.Range("N122").Value = -44895.499
$ws.Range("H126").Value = 3538.3845
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 6999.6665
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 20998.9995
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -25938.9995
$ws.Range("H132").Value = 6581341.5
$ws.Range("I132").Value = 6946810.5
$ws.Range("K132").Value = 20840431.5
$ws.Range("M132").Value = -20837901.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3919.889
$ws.Range("I7").Value = 3857.6
$ws.Range("J7").Value = 3997.75
$ws.Range("K7").Value = 3857.6
$ws.Range("L7").Value = 3997.75
$ws.Range("M7").Value = -3745.6
$ws.Range("N7").Value = -4221.75
$ws.Range("H16").Value = 2147.818
$ws.Range("J16").Value = 2153
$ws.Range("L16").Value = 2153
$ws.Range("N16").Value = -2493
$ws.Range("H22").Value = 1434.6957
$ws.Range("I22").Value = 1605.7646
$ws.Range("K22").Value = 1605.7646
$ws.Range("M22").Value = -1310.7646
$ws.Range("H27").Value = 1434.6957
$ws.Range("I27").Value = 1605.7646
$ws.Range("K27").Value = 1605.7646
$ws.Range("M27").Value = -1498.7646
$ws.Range("H34").Value = 19249.5
$ws.Range("I34").Value = 19249.5
$ws.Range("K34").Value = 19249.5
$ws.Range("M34").Value = -19046.5
$ws.Range("H38").Value = 34999
$ws.Range("J38").Value = 34999
$ws.Range("L38").Value = 34999
$ws.Range("N38").Value = -35819
$ws.Range("H40").Value = 3497.8
$ws.Range("I40").Value = 3123
$ws.Range("J40").Value = 4997
$ws.Range("K40").Value = 3123
$ws.Range("L40").Value = 4997
$ws.Range("M40").Value = -2987
$ws.Range("N40").Value = -5269
$ws.Range("H46").Value = 1251.8667
$ws.Range("I46").Value = 1655.8572
$ws.Range("K46").Value = 1655.8572
$ws.Range("M46").Value = -1467.8572
$ws.Range("H55").Value = 529.1818
$ws.Range("J55").Value = 758.75
$ws.Range("L55").Value = 758.75
$ws.Range("N55").Value = -1104.75
$ws.Range("H82").Value = 1464.4667
$ws.Range("I82").Value = 1578.2273
$ws.Range("J82").Value = 1151.625
$ws.Range("K82").Value = 1578.2273
$ws.Range("L82").Value = 1151.625
$ws.Range("M82").Value = -1217.2273
$ws.Range("N82").Value = -1873.625
$ws.Range("H85").Value = 1464.4667
$ws.Range("I85").Value = 1578.2273
$ws.Range("J85").Value = 1151.625
$ws.Range("K85").Value = 1578.2273
$ws.Range("L85").Value = 1151.625
$ws.Range("M85").Value = -330.2273
$ws.Range("N85").Value = -3647.625
$ws.Range("H93").Value = 1919.3889
$ws.Range("I93").Value = 1000.6429
$ws.Range("K93").Value = 1000.6429
$ws.Range("M93").Value = 247.3570999999999
$ws.Range("H122").Value = 11459.389
$ws.Range("I122").Value = 10580.692
$ws.Range("J122").Value = 13744
$ws.Range("K122").Value = 31742.076
$ws.Range("L122").Value = 41232
$ws.Range("M122").Value = -29292.076
$ws.Range("N122").Value = -46132
$ws.Range("H126").Value = 3919.889
$ws.Range("I126").Value = 3857.6
$ws.Range("J126").Value = 3997.75
$ws.Range("K126").Value = 11572.8
$ws.Range("L126").Value = 11993.25
$ws.Range("M126").Value = -9102.799999999999
$ws.Range("N126").Value = -16933.25
$ws.Range("H132").Value = 15486812
$ws.Range("I132").Value = 17145728
$ws.Range("J132").Value = 3589.6667
$ws.Range("K132").Value = 51437184
$ws.Range("L132").Value = 10769.0001
$ws.Range("M132").Value = -51434654
$ws.Range("N132").Value = -15829.0001
$ws.Range("H136").Value = 1242.65
$ws.Range("I136").Value = 1275.0588
$ws.Range("J136").Value = 1059
$ws.Range("K136").Value = 3825.1764
$ws.Range("L136").Value = 3177
$ws.Range("M136").Value = -1275.1764
$ws.Range("N136").Value = -8277

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 44604.652
$ws.Range("I81").Value = 46404.863
$ws.Range("K81").Value = 92809.726
$ws.Range("M81").Value = -91748.726
$ws.Range("H84").Value = 44604.652
$ws.Range("I84").Value = 46404.863
$ws.Range("K84").Value = 464048.63
$ws.Range("M84").Value = -458744.63
$ws.Range("H100").Value = 1668.2307
$ws.Range("I100").Value = 1615.3636
$ws.Range("J100").Value = 1959
$ws.Range("K100").Value = 3230.7272
$ws.Range("L100").Value = 3918
$ws.Range("M100").Value = -2689.7272
$ws.Range("N100").Value = -5000
$ws.Range("H122").Value = 1641.3334
$ws.Range("I122").Value = 1641.3334
$ws.Range("K122").Value = 4924.0002
$ws.Range("M122").Value = -2474.0002
$ws.Range("H126").Value = 3628.6
$ws.Range("I126").Value = 3628.6
$ws.Range("K126").Value = 10885.8
$ws.Range("M126").Value = -8415.799999999999
$ws.Range("H132").Value = 8068221
$ws.Range("I132").Value = 10418749
$ws.Range("J132").Value = 9268.571
$ws.Range("K132").Value = 31256247
$ws.Range("L132").Value = 27805.713
$ws.Range("M132").Value = -31253717
$ws.Range("N132").Value = -32865.713
$ws.Range("H136").Value = 13159362
$ws.Range("I136").Value = 14707314
$ws.Range("J136").Value = 1773
$ws.Range("K136").Value = 44121942
$ws.Range("L136").Value = 5319
$ws.Range("M136").Value = -44119392
$ws.Range("N136").Value = -10419
